$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '93.572.63'
$ws.Cells.Item(2, 5).Value = '  +2.02%  '

$ws.Cells.Item(3, 4).Value = '3.103.48'
$ws.Cells.Item(3, 5).Value = '  -0.78%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '237.78'
$ws.Cells.Item(5, 5).Value = '  -3.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '612.67'
$ws.Cells.Item(6, 5).Value = '  -0.95%  '

$ws.Cells.Item(7, 5).Value = '  +2.17%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.390'
$ws.Cells.Item(8, 5).Value = '  +1.70%  '

$ws.Cells.Item(9, 5).Value = '  -0.04%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.827'
$ws.Cells.Item(10, 5).Value = '  +12.06%  '

$ws.Cells.Item(11, 4).Value = '3.102.18'
$ws.Cells.Item(11, 5).Value = '  -0.76%  '

$ws.Cells.Item(12, 5).Value = '  -3.43%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.0000244'
$ws.Cells.Item(13, 5).Value = '  -2.66%  '

$ws.Cells.Item(14, 2).Value = 'Avalanche'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '34.90'
$ws.Cells.Item(14, 5).Value = '  +0.07%  '

$ws.Cells.Item(15, 2).Value = 'WrappedBTC'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(15, 4).Value = '93.255.78'
$ws.Cells.Item(15, 5).Value = '  +1.79%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '5.42'
$ws.Cells.Item(16, 5).Value = '  -3.08%  '

$ws.Cells.Item(17, 4).Value = '3.675.88'
$ws.Cells.Item(17, 5).Value = '  -0.79%  '

$ws.Cells.Item(18, 4).Value = '3.103.50'
$ws.Cells.Item(18, 5).Value = '  -1.03%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '3.66'
$ws.Cells.Item(19, 5).Value = '  -1.09%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '14.71'
$ws.Cells.Item(20, 5).Value = '  -1.77%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.03'
$ws.Cells.Item(21, 5).Value = '  +2.62%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '442.03'
$ws.Cells.Item(22, 5).Value = '  -1.65%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.0000199'
$ws.Cells.Item(23, 5).Value = '  -0.69%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.01'
$ws.Cells.Item(24, 5).Value = '  -5.09%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '8.20'
$ws.Cells.Item(25, 5).Value = '  +4.12%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '5.68'
$ws.Cells.Item(26, 5).Value = '  -4.00%  '

$ws.Cells.Item(27, 5).Value = '  +8.36%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '86.00'
$ws.Cells.Item(28, 5).Value = '  -2.32%  '

$ws.Cells.Item(29, 2).Value = 'WrappedeETH'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(29, 4).Value = '3.270.32'
$ws.Cells.Item(29, 5).Value = '  -0.62%  '

$ws.Cells.Item(30, 2).Value = 'Dai'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.999'
$ws.Cells.Item(30, 5).Value = '  -0.05%  '

$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.252'
$ws.Cells.Item(31, 5).Value = '  +6.64%  '

$ws.Cells.Item(32, 2).Value = 'Cronos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.182'
$ws.Cells.Item(32, 5).Value = '  +8.66%  '

$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.124'
$ws.Cells.Item(33, 5).Value = '  -16.14%  '

$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '9.21'
$ws.Cells.Item(34, 5).Value = '  -1.89%  '

$ws.Cells.Item(35, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +52.97%  '

$ws.Cells.Item(36, 2).Value = 'RenderToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '7.95'
$ws.Cells.Item(36, 5).Value = '  +0.23%  '

$ws.Cells.Item(37, 2).Value = 'Kaspa'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.159'
$ws.Cells.Item(37, 5).Value = '  -9.51%  '

$ws.Cells.Item(38, 2).Value = 'EthereumClassic'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '25.95'
$ws.Cells.Item(38, 5).Value = '  -1.41%  '

$ws.Cells.Item(39, 2).Value = 'MantraDAO'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.94'
$ws.Cells.Item(39, 5).Value = '  -6.87%  '

$ws.Cells.Item(40, 2).Value = 'PancakeSwap'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.90'
$ws.Cells.Item(40, 5).Value = '  -1.26%  '

$ws.Cells.Item(41, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.453'
$ws.Cells.Item(41, 5).Value = '  +2.01%  '

$ws.Cells.Item(42, 2).Value = 'WhiteBITCoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '24.00'
$ws.Cells.Item(42, 5).Value = '  +8.25%  '

$ws.Cells.Item(43, 2).Value = 'Bittensor'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '475.63'
$ws.Cells.Item(43, 5).Value = '  -3.57%  '

$ws.Cells.Item(44, 2).Value = 'Fetch.AI'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '1.28'
$ws.Cells.Item(44, 5).Value = '  -1.87%  '

$ws.Cells.Item(45, 2).Value = 'dogwifhat'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.28'
$ws.Cells.Item(45, 5).Value = '  -3.66%  '

$ws.Cells.Item(46, 2).Value = 'USDe'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.00'
$ws.Cells.Item(46, 5).Value = '  +0.01%  '

$ws.Cells.Item(47, 2).Value = 'Monero'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '159.76'
$ws.Cells.Item(47, 5).Value = '  +0.53%  '

$ws.Cells.Item(48, 2).Value = 'ARBITRUM'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.702'
$ws.Cells.Item(48, 5).Value = '  -1.13%  '

$ws.Cells.Item(49, 2).Value = 'Stacks'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.86'
$ws.Cells.Item(49, 5).Value = '  -3.00%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.32'
$ws.Cells.Item(50, 5).Value = '  -4.08%  '

$ws.Cells.Item(51, 2).Value = 'OKB'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '43.79'
$ws.Cells.Item(51, 5).Value = '  -0.59%  '
